$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")
$ws.Range("A10").Value = "ayvid"
$ws.Range("A10").Select()
